$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.90199178708178
$ws.Range("C2").Value = 8.604681879602609
$ws.Range("D2").Value = 3.822763447647317
$ws.Range("F2").Value = 19.63438735485051
$ws.Range("G2").Value = 3.594575244415813
$ws.Range("I2").Value = 16.85341872438862
$ws.Range("M2").Value = 19.96750138311285
$ws.Range("O2").Value = 17.450539578481
$ws.Range("B3").Value = 10.30525769918317
$ws.Range("C3").Value = 8.228218988018671
$ws.Range("D3").Value = 3.758950552275222
$ws.Range("F3").Value = 19.61611033640073
$ws.Range("G3").Value = 3.59654367381797
$ws.Range("I3").Value = 16.97372054518728
$ws.Range("M3").Value = 19.36066237387102
$ws.Range("O3").Value = 17.50402461552709
$ws.Range("B4").Value = 9.920151481090317
$ws.Range("C4").Value = 7.986671605146648
$ws.Range("D4").Value = 3.71881348091502
$ws.Range("F4").Value = 19.61255712549032
$ws.Range("G4").Value = 3.59781563645894
$ws.Range("I4").Value = 17.0524799283568
$ws.Range("M4").Value = 18.98641733841324
$ws.Range("O4").Value = 17.54268122214089
$ws.Range("B5").Value = 9.758635869823667
$ws.Range("C5").Value = 7.885716823440275
$ws.Range("D5").Value = 3.702229425991182
$ws.Range("F5").Value = 19.6130353703958
$ws.Range("G5").Value = 3.598349952554913
$ws.Range("I5").Value = 17.08580235709037
$ws.Range("M5").Value = 18.83375315191545
$ws.Range("O5").Value = 17.55988866564775
$ws.Range("B6").Value = 9.73154369489805
$ws.Range("C6").Value = 7.86880399590877
$ws.Range("D6").Value = 3.699462266238967
$ws.Range("F6").Value = 19.61323101242765
$ws.Range("G6").Value = 3.598439642131353
$ws.Range("I6").Value = 17.09140955852994
$ws.Range("M6").Value = 18.80840170126369
$ws.Range("O6").Value = 17.56283355552339
$ws.Range("B7").Value = 9.917991595339611
$ws.Range("C7").Value = 7.985320173926053
$ws.Range("D7").Value = 3.718590728415724
$ws.Range("F7").Value = 19.61255578064717
$ws.Range("G7").Value = 3.597822777654121
$ws.Range("I7").Value = 17.05292436090307
$ws.Range("M7").Value = 18.98435873952126
$ws.Range("O7").Value = 17.54290740923208
$ws.Range("B8").Value = 10.70018714768957
$ws.Range("C8").Value = 8.47708665108148
$ws.Range("D8").Value = 3.80096643151218
$ws.Range("F8").Value = 19.62649309114714
$ws.Range("G8").Value = 3.595240845372339
$ws.Range("I8").Value = 16.8938811819364
$ws.Range("M8").Value = 19.75875435168649
$ws.Range("O8").Value = 17.46776928360153
$ws.Range("B9").Value = 12.08137856926186
$ws.Range("C9").Value = 9.355550411847007
$ws.Range("D9").Value = 3.954401089951328
$ws.Range("F9").Value = 19.71466801954407
$ws.Range("G9").Value = 3.590677787090139
$ws.Range("I9").Value = 16.62097273869357
$ws.Range("M9").Value = 21.25374461567991
$ws.Range("O9").Value = 17.36691398763468
$ws.Range("B10").Value = 12.99868745165368
$ws.Range("C10").Value = 9.944873346085297
$ws.Range("D10").Value = 4.061530129337694
$ws.Range("F10").Value = 19.81638306162392
$ws.Range("G10").Value = 3.587626727834858
$ws.Range("I10").Value = 16.44443371926457
$ws.Range("M10").Value = 22.32474643394873
$ws.Range("O10").Value = 17.32158606168034
$ws.Range("B11").Value = 13.40817152857923
$ws.Range("C11").Value = 10.20020434636768
$ws.Range("D11").Value = 4.108921643223519
$ws.Range("F11").Value = 19.87059267935864
$ws.Range("G11").Value = 3.586303434060399
$ws.Range("I11").Value = 16.36937341199685
$ws.Range("M11").Value = 22.80363594590214
$ws.Range("O11").Value = 17.30729213960942
$ws.Range("B12").Value = 13.56187661535275
$ws.Range("C12").Value = 10.29501882333596
$ws.Range("D12").Value = 4.1266646060467
$ws.Range("F12").Value = 19.89225141932937
$ws.Range("G12").Value = 3.585811577149039
$ws.Range("I12").Value = 16.34170907509014
$ws.Range("M12").Value = 22.98360810374066
$ws.Range("O12").Value = 17.30279458720241
$ws.Range("B13").Value = 13.52893492003474
$ws.Range("C13").Value = 10.27468272795148
$ws.Range("D13").Value = 4.1228525337312
$ws.Range("F13").Value = 19.88753672973489
$ws.Range("G13").Value = 3.585917096933429
$ws.Range("I13").Value = 16.34763323503786
$ws.Range("M13").Value = 22.94491157350491
$ws.Range("O13").Value = 17.30372242826686
$ws.Range("B14").Value = 13.42088968925759
$ws.Range("C14").Value = 10.20804254784242
$ws.Range("D14").Value = 4.110385495641307
$ws.Range("F14").Value = 19.87235196570273
$ws.Range("G14").Value = 3.586262783687615
$ws.Range("I14").Value = 16.367082203871
$ws.Range("M14").Value = 22.81847083089093
$ws.Range("O14").Value = 17.30690375087164
$ws.Range("B15").Value = 13.35423613331617
$ws.Range("C15").Value = 10.16697842467115
$ws.Range("D15").Value = 4.102722323144636
$ws.Range("F15").Value = 19.86319775265941
$ws.Range("G15").Value = 3.586475729372731
$ws.Range("I15").Value = 16.37909429393575
$ws.Range("M15").Value = 22.74083847492051
$ws.Range("O15").Value = 17.30897174867511
$ws.Range("B16").Value = 12.9723974949976
$ws.Range("C16").Value = 9.927926650448189
$ws.Range("D16").Value = 4.058405092463802
$ws.Range("F16").Value = 19.81299920984657
$ws.Range("G16").Value = 3.587714504752214
$ws.Range("I16").Value = 16.44944509039632
$ws.Range("M16").Value = 22.29326700486855
$ws.Range("O16").Value = 17.32264799029565
$ws.Range("B17").Value = 12.73956079972242
$ws.Range("C17").Value = 9.777978791666072
$ws.Range("D17").Value = 4.030866739874475
$ws.Range("F17").Value = 19.78423037179231
$ws.Range("G17").Value = 3.588490975220588
$ws.Range("I17").Value = 16.49395048876018
$ws.Range("M17").Value = 22.01643501707193
$ws.Range("O17").Value = 17.33266233984829
$ws.Range("B18").Value = 12.60359245114301
$ws.Range("C18").Value = 9.690534083604495
$ws.Range("D18").Value = 4.014901661880395
$ws.Range("F18").Value = 19.76843134569217
$ws.Range("G18").Value = 3.58894366852376
$ws.Range("I18").Value = 16.52004250837544
$ws.Range("M18").Value = 21.85643542178042
$ws.Range("O18").Value = 17.33901746871044
$ws.Range("B19").Value = 12.55720570192442
$ws.Range("C19").Value = 9.660722260505429
$ws.Range("D19").Value = 4.009474884199571
$ws.Range("F19").Value = 19.76321083875672
$ws.Range("G19").Value = 3.58909798986747
$ws.Range("I19").Value = 16.52896146991567
$ws.Range("M19").Value = 21.80213521240361
$ws.Range("O19").Value = 17.34127124085796
$ws.Range("B20").Value = 12.7645588156597
$ws.Range("C20").Value = 9.794065364030178
$ws.Range("D20").Value = 4.033811334485985
$ws.Range("F20").Value = 19.7872155150631
$ws.Range("G20").Value = 3.588407688893889
$ws.Range("I20").Value = 16.48916168454792
$ws.Range("M20").Value = 22.04598561457517
$ws.Range("O20").Value = 17.33153465797574
$ws.Range("B21").Value = 13.45272370993435
$ws.Range("C21").Value = 10.22766752174013
$ws.Range("D21").Value = 4.114052959309922
$ws.Range("F21").Value = 19.87678150911103
$ws.Range("G21").Value = 3.586160996584371
$ws.Range("I21").Value = 16.36134892281676
$ws.Range("M21").Value = 22.85564812684195
$ws.Range("O21").Value = 17.30594443874687
$ws.Range("B22").Value = 13.89335634962653
$ws.Range("C22").Value = 10.50011569933144
$ws.Range("D22").Value = 4.165306934726969
$ws.Range("F22").Value = 19.94190221546919
$ws.Range("G22").Value = 3.584746521084343
$ws.Range("I22").Value = 16.28224440716433
$ws.Range("M22").Value = 23.37673724923724
$ws.Range("O22").Value = 17.29455660108878
$ws.Range("B23").Value = 13.66011829595604
$ws.Range("C23").Value = 10.35571700006107
$ws.Range("D23").Value = 4.138063690094262
$ws.Range("F23").Value = 19.90654779476408
$ws.Range("G23").Value = 3.585496541043168
$ws.Range("I23").Value = 16.32405716845561
$ws.Range("M23").Value = 23.09941477102058
$ws.Range("O23").Value = 17.30014451784481
$ws.Range("B24").Value = 12.75326377377807
$ws.Range("C24").Value = 9.78679647670041
$ws.Range("D24").Value = 4.032480496474093
$ws.Range("F24").Value = 19.78586362430737
$ws.Range("G24").Value = 3.588445323046057
$ws.Range("I24").Value = 16.49132512947184
$ws.Range("M24").Value = 22.03262840840087
$ws.Range("O24").Value = 17.33204262153586
$ws.Range("B25").Value = 11.72458996854365
$ws.Range("C25").Value = 9.127522596027021
$ws.Range("D25").Value = 3.913830716589279
$ws.Range("F25").Value = 19.6843077917507
$ws.Range("G25").Value = 3.591859033998849
$ws.Range("I25").Value = 16.69060845779786
$ws.Range("M25").Value = 20.85324350144308
$ws.Range("O25").Value = 17.38917439015208
